{"js": "// Update the division-problem worksheet: replace each \"old\" expression\n// with its corresponding \"new\" expression throughout the document body.\n// Every \"old\" string is a unique substring in the document, so a single\n// search + replace per pair is safe.\n\nconst replacements = [\n  [\"62\u00f78=\", \"45\u00f79=\"],\n  [\"64\u00f73=\", \"60\u00f74=\"],\n  [\"22\u00f73=\", \"45\u00f75=\"],\n  [\"38\u00f75=\", \"58\u00f74=\"],\n  [\"33\u00f79=\", \"48\u00f72=\"],\n  [\"24\u00f79=\", \"26\u00f79=\"],\n  [\"94\u00f76=\", \"73\u00f75=\"],\n  [\"48\u00f74=\", \"84\u00f73=\"],\n  [\"56\u00f79=\", \"30\u00f79=\"],\n  [\"42\u00f75=\", \"43\u00f74=\"],\n  [\"36\u00f72=\", \"67\u00f74=\"],\n  [\"90\u00f79=\", \"33\u00f72=\"],\n  [\"18\u00f72=\", \"78\u00f76=\"],\n  [\"47\u00f79=\", \"64\u00f76=\"],\n  [\"80\u00f79=\", \"93\u00f76=\"],\n  [\"83\u00f75=\", \"88\u00f79=\"],\n  [\"40\u00f78=\", \"68\u00f76=\"],\n  [\"27\u00f76=\", \"58\u00f73=\"],\n  [\"12\u00f79=\", \"76\u00f78=\"],\n  [\"35\u00f72=\", \"35\u00f79=\"],\n  [\"69\u00f74=\", \"10\u00f79=\"],\n  [\"33\u00f77=\", \"86\u00f75=\"],\n  [\"75\u00f76=\", \"83\u00f78=\"],\n  [\"61\u00f74=\", \"18\u00f72=\"],\n  [\"82\u00f74=\", \"90\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem worksheet: replace each \"old\" expression\n# with its corresponding \"new\" expression throughout the document body.\n# Every \"old\" string is a unique substring in the document, so a scoped\n# Find/Replace (wdReplaceAll) per pair is safe and idempotent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"62\u00f78=\"; New = \"45\u00f79=\" },\n    @{ Old = \"64\u00f73=\"; New = \"60\u00f74=\" },\n    @{ Old = \"22\u00f73=\"; New = \"45\u00f75=\" },\n    @{ Old = \"38\u00f75=\"; New = \"58\u00f74=\" },\n    @{ Old = \"33\u00f79=\"; New = \"48\u00f72=\" },\n    @{ Old = \"24\u00f79=\"; New = \"26\u00f79=\" },\n    @{ Old = \"94\u00f76=\"; New = \"73\u00f75=\" },\n    @{ Old = \"48\u00f74=\"; New = \"84\u00f73=\" },\n    @{ Old = \"56\u00f79=\"; New = \"30\u00f79=\" },\n    @{ Old = \"42\u00f75=\"; New = \"43\u00f74=\" },\n    @{ Old = \"36\u00f72=\"; New = \"67\u00f74=\" },\n    @{ Old = \"90\u00f79=\"; New = \"33\u00f72=\" },\n    @{ Old = \"18\u00f72=\"; New = \"78\u00f76=\" },\n    @{ Old = \"47\u00f79=\"; New = \"64\u00f76=\" },\n    @{ Old = \"80\u00f79=\"; New = \"93\u00f76=\" },\n    @{ Old = \"83\u00f75=\"; New = \"88\u00f79=\" },\n    @{ Old = \"40\u00f78=\"; New = \"68\u00f76=\" },\n    @{ Old = \"27\u00f76=\"; New = \"58\u00f73=\" },\n    @{ Old = \"12\u00f79=\"; New = \"76\u00f78=\" },\n    @{ Old = \"35\u00f72=\"; New = \"35\u00f79=\" },\n    @{ Old = \"69\u00f74=\"; New = \"10\u00f79=\" },\n    @{ Old = \"33\u00f77=\"; New = \"86\u00f75=\" },\n    @{ Old = \"75\u00f76=\"; New = \"83\u00f78=\" },\n    @{ Old = \"61\u00f74=\"; New = \"18\u00f72=\" },\n    @{ Old = \"82\u00f74=\"; New = \"90\u00f75=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
